# PALI_rNPV_Model.xlsx -- "Update future dilution estimate"
#
# The Phase 3 funding-event assumption on the "Diluted Shares" sheet is
# revised: the projected ~$60M raise is now assumed to price at ~$6.00/share
# (was ~$6.50/share), which implies 10,000,000 new shares instead of
# 9,230,769. The footnote text explaining the assumption is updated to
# match. All downstream formulas (fully diluted share count, fair value
# per share, sensitivity tables, etc.) recalculate automatically from
# this change.
#
# Also, the exercise price recorded for the already fully-exercised Oct
# 2025 pre-funded warrants is corrected to a nominal amount.

$wb = $excel.ActiveWorkbook

$wsDiluted = $wb.Worksheets.Item("Diluted Shares")
$wsModel   = $wb.Worksheets.Item("rNPV Model")
$wsSens    = $wb.Worksheets.Item("Sensitivity")

# --- Diluted Shares sheet -------------------------------------------------

# Row 8: "Oct 2025 Pre-Funded Warrants (fully exercised)" — exercise price
# is effectively nominal (pre-funded warrants are exercised for a token
# amount), correct it from 0.7 down to 0.0001.
$wsDiluted.Range("C8").Value = 0.0001

# Row 21: "Future equity raise for Phase 3 funding (est.)" — revise the
# assumed raise price from ~$6.50/share to ~$6.00/share, which changes the
# implied share count from 9,230,769 to 10,000,000, and update the note.
$wsDiluted.Range("B21").Value = 10000000
$wsDiluted.Range("C21").Value = 6
$wsDiluted.Range("D21").Value = "Assumption: ~`$60M raise at ~`$6.00/share Phase 3 funding event (2028)"

# --- Mirror the author's on-screen navigation/selection at save time ------

$wsModel.Activate()
$wsModel.Range("B21").Select()

$wsSens.Activate()
$wsSens.Range("H21").Select()

$wsDiluted.Activate()
$wsDiluted.Range("A25").Select()
